# Auto-generated edit script: update Gilgamesh_Profits (per-sheet leve profit data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3145.115
$ws.Range("I15").Value = 3145.115
$ws.Range("K15").Value = 9435.344999999999
$ws.Range("M15").Value = -9266.344999999999
$ws.Range("H18").Value = 83337864
$ws.Range("I18").Value = 4813.9
$ws.Range("K18").Value = 4813.9
$ws.Range("M18").Value = -4529.9
$ws.Range("H40").Value = 4444.8
$ws.Range("I40").Value = 4112.4165
$ws.Range("K40").Value = 4112.4165
$ws.Range("M40").Value = -3937.4165
$ws.Range("H132").Value = 5314.2974
$ws.Range("I132").Value = 5538
$ws.Range("J132").Value = 1399.5
$ws.Range("K132").Value = 16614
$ws.Range("L132").Value = 4198.5
$ws.Range("M132").Value = -14084
$ws.Range("N132").Value = -9258.5
$ws.Range("H138").Value = 280945.1
$ws.Range("I138").Value = 2997.6758
$ws.Range("J138").Value = 474983.88
$ws.Range("K138").Value = 8993.027399999999
$ws.Range("L138").Value = 1424951.64
$ws.Range("M138").Value = -3853.027399999999
$ws.Range("N138").Value = -1435231.64
$ws.Range("H141").Value = 4012.818
$ws.Range("I141").Value = 3529.5
$ws.Range("K141").Value = 10588.5
$ws.Range("M141").Value = -5408.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4193.727
$ws.Range("I32").Value = 3685.0327
$ws.Range("K32").Value = 3685.0327
$ws.Range("M32").Value = -3398.0327
$ws.Range("H45").Value = 19828.566
$ws.Range("I45").Value = 24371.592
$ws.Range("J45").Value = 7335.25
$ws.Range("K45").Value = 24371.592
$ws.Range("L45").Value = 7335.25
$ws.Range("M45").Value = -23994.592
$ws.Range("N45").Value = -8089.25
$ws.Range("H61").Value = 4917.8
$ws.Range("I61").Value = 2945.0833
$ws.Range("K61").Value = 2945.0833
$ws.Range("M61").Value = -2733.0833
$ws.Range("H88").Value = 7361.2
$ws.Range("I88").Value = 4998
$ws.Range("J88").Value = 8936.666999999999
$ws.Range("K88").Value = 4998
$ws.Range("L88").Value = 8936.666999999999
$ws.Range("M88").Value = -4592
$ws.Range("N88").Value = -9748.666999999999
$ws.Range("H91").Value = 7361.2
$ws.Range("I91").Value = 4998
$ws.Range("J91").Value = 8936.666999999999
$ws.Range("K91").Value = 4998
$ws.Range("L91").Value = 8936.666999999999
$ws.Range("M91").Value = -3594
$ws.Range("N91").Value = -11744.667
$ws.Range("H132").Value = 3621.348
$ws.Range("I132").Value = 2582
$ws.Range("K132").Value = 7746
$ws.Range("M132").Value = -5216
$ws.Range("H136").Value = 4917.8
$ws.Range("I136").Value = 2945.0833
$ws.Range("K136").Value = 8835.249899999999
$ws.Range("M136").Value = -6285.249899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 55555824
$ws.Range("I94").Value = 62500210
$ws.Range("J94").Value = 749.5
$ws.Range("K94").Value = 62500210
$ws.Range("L94").Value = 749.5
$ws.Range("M94").Value = -62499759
$ws.Range("N94").Value = -1651.5
$ws.Range("H105").Value = 13687050
$ws.Range("I105").Value = 836000.8
$ws.Range("K105").Value = 836000.8
$ws.Range("M105").Value = -834253.8
$ws.Range("H134").Value = 2932.5557
$ws.Range("I134").Value = 2214.3076
$ws.Range("K134").Value = 6642.9228
$ws.Range("M134").Value = -4107.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 20022800
$ws.Range("I62").Value = 33341998
$ws.Range("K62").Value = 33341998
$ws.Range("M62").Value = -33341374
$ws.Range("H65").Value = 20022800
$ws.Range("I65").Value = 33341998
$ws.Range("K65").Value = 166709990
$ws.Range("M65").Value = -166706870
$ws.Range("H86").Value = 1164852.6
$ws.Range("I86").Value = 3235.1538
$ws.Range("J86").Value = 3322142.2
$ws.Range("K86").Value = 3235.1538
$ws.Range("L86").Value = 3322142.2
$ws.Range("M86").Value = -2112.1538
$ws.Range("N86").Value = -3324388.2
$ws.Range("H89").Value = 1164852.6
$ws.Range("I89").Value = 3235.1538
$ws.Range("J89").Value = 3322142.2
$ws.Range("K89").Value = 16175.769
$ws.Range("L89").Value = 16610711
$ws.Range("M89").Value = -10559.769
$ws.Range("N89").Value = -16621943
$ws.Range("H134").Value = 6245.3335
$ws.Range("I134").Value = 6027.7334
$ws.Range("K134").Value = 18083.2002
$ws.Range("M134").Value = -15548.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1250
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = -1040
$ws.Range("N51").Value = -6920
$ws.Range("H127").Value = 989.4286
$ws.Range("J127").Value = 989.4286
$ws.Range("L127").Value = 2968.2858
$ws.Range("N127").Value = -12888.2858
$ws.Range("H136").Value = 739.8333
$ws.Range("I136").Value = 739.8333
$ws.Range("K136").Value = 2219.4999
$ws.Range("M136").Value = 2880.5001
$ws.Range("H139").Value = 3069.4707
$ws.Range("I139").Value = 2247.625
$ws.Range("K139").Value = 6742.875
$ws.Range("M139").Value = -1602.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3704.3635
$ws.Range("I132").Value = 3724.6667
$ws.Range("K132").Value = 11174.0001
$ws.Range("M132").Value = -8644.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1162.2727
$ws.Range("I16").Value = 1112.5555
$ws.Range("J16").Value = 1386
$ws.Range("K16").Value = 1112.5555
$ws.Range("L16").Value = 1386
$ws.Range("M16").Value = -942.5554999999999
$ws.Range("N16").Value = -1726
$ws.Range("H40").Value = 24840.959
$ws.Range("I40").Value = 27755.477
$ws.Range("K40").Value = 27755.477
$ws.Range("M40").Value = -27619.477
$ws.Range("H46").Value = 1420.8889
$ws.Range("I46").Value = 1006.3333
$ws.Range("K46").Value = 1006.3333
$ws.Range("M46").Value = -818.3333
$ws.Range("H61").Value = 1587.931
$ws.Range("I61").Value = 1473.7693
$ws.Range("J61").Value = 2577.3333
$ws.Range("K61").Value = 1473.7693
$ws.Range("L61").Value = 2577.3333
$ws.Range("M61").Value = -1271.7693
$ws.Range("N61").Value = -2981.3333
$ws.Range("H68").Value = 9000
$ws.Range("I68").Value = 9000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -8251
$ws.Range("H71").Value = 9000
$ws.Range("I71").Value = 9000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 45000
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -41256
$ws.Range("H93").Value = 1161.3125
$ws.Range("I93").Value = 1129.5385
$ws.Range("K93").Value = 1129.5385
$ws.Range("M93").Value = 118.4614999999999
$ws.Range("H113").Value = 1587.931
$ws.Range("I113").Value = 1473.7693
$ws.Range("J113").Value = 2577.3333
$ws.Range("K113").Value = 1473.7693
$ws.Range("L113").Value = 2577.3333
$ws.Range("M113").Value = 696.2307000000001
$ws.Range("N113").Value = -6917.3333
$ws.Range("H122").Value = 2559.6
$ws.Range("I122").Value = 1866.3334
$ws.Range("K122").Value = 5599.0002
$ws.Range("M122").Value = -3149.0002
$ws.Range("H132").Value = 8242.375
$ws.Range("I132").Value = 6137.4
$ws.Range("K132").Value = 18412.2
$ws.Range("M132").Value = -15882.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 22628
$ws.Range("J54").Value = 49326
$ws.Range("L54").Value = 49326
$ws.Range("N54").Value = -50366
$ws.Range("H132").Value = 10419865
$ws.Range("I132").Value = 13336356
$ws.Range("J132").Value = 3828.8572
$ws.Range("K132").Value = 40009068
$ws.Range("L132").Value = 11486.5716
$ws.Range("M132").Value = -40006538
$ws.Range("N132").Value = -16546.5716
$ws.Range("H136").Value = 25002418
$ws.Range("I136").Value = 29412916
$ws.Range("J136").Value = 9595
$ws.Range("K136").Value = 88238748
$ws.Range("L136").Value = 28785
$ws.Range("M136").Value = -88236198
$ws.Range("N136").Value = -33885
